# Update Work Week and Social Spending
# (commit message is generic; the actual edit refreshes the GDP per Capita
# series for Kyrgyzstan on the "Data" sheet and extends it from 2010 through
# 2016.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New "Data" (column E) values keyed by Year (column D). Years 1974-1979 are
# intentionally omitted below - their cells were (and remain) blank, so they
# are left untouched. Row number = Year - 1971.
$values = [ordered]@{
    1973 = "5941"
    1980 = "5848"
    1981 = "5789"
    1982 = "5726"
    1983 = "5992"
    1984 = "5952"
    1985 = "5616"
    1986 = "5413"
    1987 = "5241"
    1988 = "5679"
    1989 = "5719"
    1990 = "5742"
    1991 = "5161.90411604431"
    1992 = "4362.11265416902"
    1993 = "3764.94211236256"
    1994 = "2993.82781936695"
    1995 = "2776.03031021925"
    1996 = "2901.55419395731"
    1997 = "3114.67495634521"
    1998 = "3105.96104991312"
    1999 = "3144.46115496837"
    2000 = "3246.38306906804"
    2001 = "3357.99300904187"
    2002 = "3299.24599974929"
    2003 = "3464.82014694552"
    2004 = "3632.1470955154"
    2005 = "3554.18565533234"
    2006 = "3594.56523758805"
    2007 = "3831.49481986061"
    2008 = "4046.72078648019"
    2009 = "4075.83393999403"
    2010 = "3977.64889256546"
    2011 = "4142"
    2012 = "4104"
    2013 = "4509"
    2014 = "4644"
    2015 = "4754"
    2016 = "4879"
}

# New rows for 2011-2016 must be appended (as rows 40-45), keeping the same
# Country Code / Country Name / Indicator values as the rest of the sheet.
$newYears = @(2011, 2012, 2013, 2014, 2015, 2016)

foreach ($year in $values.Keys) {
    $row = $year - 1971

    if ($newYears -contains $year) {
        $ws.Cells.Item($row, 1).Value = 417
        $ws.Cells.Item($row, 2).Value = "Kyrgyzstan"
        $ws.Cells.Item($row, 3).Value = "GDP per Capita"
        $ws.Cells.Item($row, 4).Value = $year
    }

    # The Data column stores these figures as text (not numbers) in the
    # original workbook, so prefix with an apostrophe to force text entry
    # the same way a user typing into Excel would.
    $ws.Cells.Item($row, 5).Value = "'" + $values[$year]
}
